{"js": "// The commit only changes the VISIBLE text in two small spots (the rest of\n// the underlying XML diff is just Word re-splitting/merging <w:r> runs and\n// removing <w:proofErr> spell-check markers around the edited words, which\n// carries no semantic/visible change):\n//   1) \"PLease\"        -> \"Please\"        (capitalization fix)\n//   2) \"please launch\" -> \"you can launch\" (wording change)\nconst body = context.document.body;\n\n// 1) Fix \"PLease\" -> \"Please\" (case-sensitive search so we hit only this\n//    exact occurrence, not the other \"Please\"/\"please\" spellings nearby).\nconst badCap = body.search(\"PLease\", { matchCase: true, matchWholeWord: false });\nbadCap.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < badCap.items.length; i++) {\n  badCap.items[i].insertText(\"Please\", \"Replace\");\n}\nawait context.sync();\n\n// 2) \"please launch\" -> \"you can launch\" in the bode-plot sentence.\nconst oldPhrase = body.search(\"please launch\", { matchCase: true, matchWholeWord: false });\noldPhrase.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < oldPhrase.items.length; i++) {\n  oldPhrase.items[i].insertText(\"you can launch\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# The commit only changes the VISIBLE text in two small spots (the rest of\n# the underlying XML diff is just Word re-splitting/merging <w:r> runs and\n# removing <w:proofErr> spell-check markers around the edited words, which\n# carries no semantic/visible change):\n#   1) \"PLease\"        -> \"Please\"        (capitalization fix)\n#   2) \"please launch\" -> \"you can launch\" (wording change)\n\n$d = $word.ActiveDocument\n\n# 1) Fix \"PLease\" -> \"Please\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\n    [ref]\"PLease\",      # FindText\n    [ref]$true,         # MatchCase\n    [ref]$false,        # MatchWholeWord\n    [ref]$false,        # MatchWildcards\n    [ref]$false,        # MatchSoundsLike\n    [ref]$false,        # MatchAllWordForms\n    [ref]$true,         # Forward\n    [ref]1,             # Wrap (wdFindContinue)\n    [ref]$false,        # Format\n    [ref]\"Please\",      # ReplaceWith\n    [ref]2              # Replace (wdReplaceAll)\n)\n\n# 2) \"please launch\" -> \"you can launch\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    [ref]\"please launch\",   # FindText\n    [ref]$true,             # MatchCase\n    [ref]$false,            # MatchWholeWord\n    [ref]$false,            # MatchWildcards\n    [ref]$false,            # MatchSoundsLike\n    [ref]$false,            # MatchAllWordForms\n    [ref]$true,             # Forward\n    [ref]1,                 # Wrap (wdFindContinue)\n    [ref]$false,            # Format\n    [ref]\"you can launch\",  # ReplaceWith\n    [ref]2                  # Replace (wdReplaceAll)\n)\n"}
